$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting from the neighboring
# header cell (G1) so it picks up the same bold/border/alignment style
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Save column values (0/1 flags), row 7 is 1, rest are 0
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 0
    10 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
